$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.28498633333333
$ws.Range("H2").Value = 36.854959
$ws.Range("I2").Value = 0.1279589698403688
$ws.Range("J2").Value = 0.1279589698403688
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.460132666666667
$ws.Range("N2").Value = 10.380398
$ws.Range("O2").Value = 0.01616897968344663
$ws.Range("P2").Value = 0.01616897968344663
$ws.Range("Q2").Value = 42.50768252152022
$ws.Range("R2").Value = 382.569142693682
$ws.Range("S2").Value = 0.002068965983663684
$ws.Range("T2").Value = 0.002068965983663684
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.28498633333333
$ws.Range("H3").Value = 36.854959
$ws.Range("I3").Value = 0.1279589698403688
$ws.Range("J3").Value = 0.1279589698403688
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 183.09167
$ws.Range("N3").Value = 549.27501
$ws.Range("O3").Value = 0.8555757185143522
$ws.Range("P3").Value = 0.8555757185143523
$ws.Range("Q3").Value = 2249.278663697177
$ws.Range("R3").Value = 20243.50797327459
$ws.Range("S3").Value = 0.1094785875615299
$ws.Range("T3").Value = 0.1094785875615299
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.28498633333333
$ws.Range("H4").Value = 36.854959
$ws.Range("I4").Value = 0.1279589698403688
$ws.Range("J4").Value = 0.1279589698403688
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 27.44640466666667
$ws.Range("N4").Value = 82.339214
$ws.Range("O4").Value = 0.1282553018022011
$ws.Range("P4").Value = 0.1282553018022011
$ws.Range("Q4").Value = 337.1787062291363
$ws.Range("R4").Value = 3034.608356062226
$ws.Range("S4").Value = 0.01641141629517525
$ws.Range("T4").Value = 0.01641141629517525
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 65.605124
$ws.Range("H5").Value = 196.815372
$ws.Range("I5").Value = 0.6833352399026945
$ws.Range("J5").Value = 0.6833352399026945
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.460132666666667
$ws.Range("N5").Value = 10.380398
$ws.Range("O5").Value = 0.01616897968344663
$ws.Range("P5").Value = 0.01616897968344663
$ws.Range("Q5").Value = 227.0024326531173
$ws.Range("R5").Value = 2043.021893878056
$ws.Range("S5").Value = 0.0110488336109698
$ws.Range("T5").Value = 0.0110488336109698
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 65.605124
$ws.Range("H6").Value = 196.815372
$ws.Range("I6").Value = 0.6833352399026945
$ws.Range("J6").Value = 0.6833352399026945
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 183.09167
$ws.Range("N6").Value = 549.27501
$ws.Range("O6").Value = 0.8555757185143522
$ws.Range("P6").Value = 0.8555757185143523
$ws.Range("Q6").Value = 12011.75171371708
$ws.Range("R6").Value = 108105.7654234537
$ws.Range("S6").Value = 0.5846450388659251
$ws.Range("T6").Value = 0.5846450388659252
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 65.605124
$ws.Range("H7").Value = 196.815372
$ws.Range("I7").Value = 0.6833352399026945
$ws.Range("J7").Value = 0.6833352399026945
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.44640466666667
$ws.Range("N7").Value = 82.339214
$ws.Range("O7").Value = 0.1282553018022011
$ws.Range("P7").Value = 0.1282553018022011
$ws.Range("Q7").Value = 1800.624781510845
$ws.Range("R7").Value = 16205.62303359761
$ws.Range("S7").Value = 0.08764136742579956
$ws.Range("T7").Value = 0.08764136742579957
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.11712033333333
$ws.Range("H8").Value = 54.351361
$ws.Range("I8").Value = 0.1887057902569366
$ws.Range("J8").Value = 0.1887057902569366
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.460132666666667
$ws.Range("N8").Value = 10.380398
$ws.Range("O8").Value = 0.01616897968344663
$ws.Range("P8").Value = 0.01616897968344663
$ws.Range("Q8").Value = 62.68763989129755
$ws.Range("R8").Value = 564.1887590216779
$ws.Range("S8").Value = 0.003051180088813149
$ws.Range("T8").Value = 0.003051180088813149
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.11712033333333
$ws.Range("H9").Value = 54.351361
$ws.Range("I9").Value = 0.1887057902569366
$ws.Range("J9").Value = 0.1887057902569366
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 183.09167
$ws.Range("N9").Value = 549.27501
$ws.Range("O9").Value = 0.8555757185143522
$ws.Range("P9").Value = 0.8555757185143523
$ws.Range("Q9").Value = 3317.093817420956
$ws.Range("R9").Value = 29853.84435678861
$ws.Range("S9").Value = 0.1614520920868972
$ws.Range("T9").Value = 0.1614520920868972
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.11712033333333
$ws.Range("H10").Value = 54.351361
$ws.Range("I10").Value = 0.1887057902569366
$ws.Range("J10").Value = 0.1887057902569366
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 27.44640466666667
$ws.Range("N10").Value = 82.339214
$ws.Range("O10").Value = 0.1282553018022011
$ws.Range("P10").Value = 0.1282553018022011
$ws.Range("Q10").Value = 497.2498160633615
$ws.Range("R10").Value = 4475.248344570254
$ws.Range("S10").Value = 0.02420251808122626
$ws.Range("T10").Value = 0.02420251808122626
